# Update imputed values in the KNN result data sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.885
$ws.Range("C7").Value = -13.163
$ws.Range("A9").Value = -21.831
$ws.Range("C12").Value = -11.555
$ws.Range("C14").Value = -12.89
$ws.Range("D15").Value = -8.348000000000003
$ws.Range("A18").Value = -22.151
$ws.Range("A20").Value = -20.558
$ws.Range("C26").Value = -13.131
$ws.Range("A27").Value = -21.761
$ws.Range("C27").Value = -13.363
$ws.Range("C29").Value = -11.857
$ws.Range("D33").Value = -7.318
$ws.Range("A35").Value = -19.732
$ws.Range("D35").Value = -7.443000000000001
$ws.Range("C37").Value = -13.151
$ws.Range("C38").Value = -13.818
$ws.Range("D38").Value = -7.825999999999999
$ws.Range("D43").Value = -7.534999999999999
$ws.Range("D44").Value = -7.616
$ws.Range("D47").Value = -7.479000000000001
$ws.Range("C51").Value = -12.405
$ws.Range("D51").Value = -7.556
$ws.Range("C52").Value = -11.361
$ws.Range("C55").Value = -13.752
$ws.Range("D57").Value = -8.038
$ws.Range("D63").Value = -7.665000000000001
$ws.Range("A69").Value = -21.862
$ws.Range("C69").Value = -11.78
$ws.Range("C70").Value = -12.579
$ws.Range("D70").Value = -7.858
$ws.Range("A76").Value = -20.306
$ws.Range("A78").Value = -19.854
$ws.Range("C81").Value = -13.422
$ws.Range("A82").Value = -22.154
$ws.Range("A83").Value = -21.758
$ws.Range("C83").Value = -12.922
$ws.Range("D88").Value = -7.816999999999998
$ws.Range("A93").Value = -21.758
$ws.Range("D99").Value = -8.103999999999999
$ws.Range("C102").Value = -13.424
